$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-21 Wednesday" "2024-02-22 Thursday"

Replace-Text "12×26=312" "51×27=1377"
Replace-Text "27×17=459" "14×11=154"
Replace-Text "73×31=2263" "99×83=8217"
Replace-Text "35×33=1155" "69×11=759"
Replace-Text "79×46=3634" "75×52=3900"

Replace-Text "61×14=854" "99×65=6435"
Replace-Text "51×70=3570" "51×24=1224"
Replace-Text "52×44=2288" "89×56=4984"
Replace-Text "41×81=3321" "88×40=3520"
Replace-Text "16×40=640" "34×39=1326"

Replace-Text "51×55=2805" "34×98=3332"
Replace-Text "96×99=9504" "19×52=988"
Replace-Text "63×40=2520" "55×93=5115"
Replace-Text "69×95=6555" "77×83=6391"
Replace-Text "26×54=1404" "34×11=374"

Replace-Text "43×44=1892" "32×88=2816"
Replace-Text "19×90=1710" "98×97=9506"
Replace-Text "96×68=6528" "64×59=3776"
Replace-Text "52×88=4576" "52×34=1768"
Replace-Text "67×42=2814" "96×46=4416"

Replace-Text "12×93=1116" "96×93=8928"
Replace-Text "16×75=1200" "48×30=1440"
Replace-Text "57×85=4845" "39×53=2067"
Replace-Text "11×39=429" "29×75=2175"
Replace-Text "43×32=1376" "63×81=5103"
